$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1274.8334
$ws.Range("J12").Value = 900
$ws.Range("L12").Value = 900
$ws.Range("N12").Value = -1240

$ws.Range("H19").Value = 1148.15
$ws.Range("I19").Value = 834.6
$ws.Range("J19").Value = 1252.6666
$ws.Range("K19").Value = 834.6
$ws.Range("L19").Value = 1252.6666
$ws.Range("M19").Value = -659.6
$ws.Range("N19").Value = -1602.6666

$ws.Range("H62").Value = 30265.166
$ws.Range("I62").Value = 43446.6
$ws.Range("J62").Value = 10880.706
$ws.Range("K62").Value = 43446.6
$ws.Range("L62").Value = 10880.706
$ws.Range("M62").Value = -42822.6
$ws.Range("N62").Value = -12128.706

$ws.Range("H65").Value = 30265.166
$ws.Range("I65").Value = 43446.6
$ws.Range("J65").Value = 10880.706
$ws.Range("K65").Value = 217233
$ws.Range("L65").Value = 54403.53
$ws.Range("M65").Value = -214113
$ws.Range("N65").Value = -60643.53

$ws.Range("H76").Value = 4456.5835
$ws.Range("I76").Value = 3622.25
$ws.Range("K76").Value = 3622.25
$ws.Range("M76").Value = -3307.25

$ws.Range("H79").Value = 4456.5835
$ws.Range("I79").Value = 3622.25
$ws.Range("K79").Value = 3622.25
$ws.Range("M79").Value = -2530.25

$ws.Range("H93").Value = 73499.5
$ws.Range("J93").Value = 73499.5
$ws.Range("L93").Value = 73499.5
$ws.Range("N93").Value = -78491.5

$ws.Range("H98").Value = 904.6111
$ws.Range("I98").Value = 863.5714
$ws.Range("K98").Value = 863.5714
$ws.Range("M98").Value = 634.4286

$ws.Range("H107").Value = 460.45456
$ws.Range("I107").Value = 465.74194
$ws.Range("K107").Value = 465.74194
$ws.Range("M107").Value = 1454.25806

$ws.Range("H112").Value = 78579.69500000001
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 78579.69500000001
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 235739.085
$ws.Range("M112").ClearContents() | Out-Null
$ws.Range("N112").Value = -237955.085

$ws.Range("H122").Value = 904.6111
$ws.Range("I122").Value = 863.5714
$ws.Range("K122").Value = 2590.7142
$ws.Range("M122").Value = -140.7142000000003

$ws.Range("H130").Value = 149995
$ws.Range("J130").Value = 149995
$ws.Range("L130").Value = 149995
$ws.Range("N130").Value = -160035

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 49995
$ws.Range("J95").Value = 49995
$ws.Range("L95").Value = 49995
$ws.Range("N95").Value = -55487

$ws.Range("H104").Value = 9999
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 9999
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 9999
$ws.Range("M104").ClearContents() | Out-Null
$ws.Range("N104").Value = -16987

$ws.Range("H112").Value = 33393.57
$ws.Range("I112").Value = 24500
$ws.Range("J112").Value = 36951
$ws.Range("K112").Value = 24500
$ws.Range("L112").Value = 36951
$ws.Range("M112").Value = -23023
$ws.Range("N112").Value = -39905

$ws.Range("H122").Value = 2772.5
$ws.Range("I122").Value = 2074.6667
$ws.Range("J122").Value = 3470.3333
$ws.Range("K122").Value = 6224.000100000001
$ws.Range("L122").Value = 10410.9999
$ws.Range("M122").Value = -3774.000100000001
$ws.Range("N122").Value = -15310.9999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 117400
$ws.Range("J20").Value = 117400
$ws.Range("L20").Value = 117400
$ws.Range("N20").Value = -117872

$ws.Range("H30").Value = 117400
$ws.Range("J30").Value = 117400
$ws.Range("L30").Value = 117400
$ws.Range("N30").Value = -117582

$ws.Range("H31").Value = 30305.621
$ws.Range("I31").Value = 39923.04
$ws.Range("K31").Value = 39923.04
$ws.Range("M31").Value = -39628.04

$ws.Range("H32").Value = 8999.875
$ws.Range("I32").Value = 15025
$ws.Range("J32").Value = 2974.75
$ws.Range("K32").Value = 15025
$ws.Range("L32").Value = 2974.75
$ws.Range("M32").Value = -14709
$ws.Range("N32").Value = -3606.75

$ws.Range("H34").Value = 30305.621
$ws.Range("I34").Value = 39923.04
$ws.Range("K34").Value = 39923.04
$ws.Range("M34").Value = -39721.04

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents() | Out-Null
$ws.Range("N41").ClearContents() | Out-Null

$ws.Range("H43").Value = 23999.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 23999.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 23999.5
$ws.Range("M43").ClearContents() | Out-Null
$ws.Range("N43").Value = -24367.5

$ws.Range("H101").Value = 23999.5
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 23999.5
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 23999.5
$ws.Range("M101").ClearContents() | Out-Null
$ws.Range("N101").Value = -30489.5

$ws.Range("H128").Value = 117400
$ws.Range("J128").Value = 117400
$ws.Range("L128").Value = 117400
$ws.Range("N128").Value = -127360

$ws.Range("H132").Value = 4924.88
$ws.Range("I132").Value = 4823.727
$ws.Range("K132").Value = 14471.181
$ws.Range("M132").Value = -11941.181

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 960.5
$ws.Range("I34").Value = 38
$ws.Range("J34").Value = 2498
$ws.Range("K34").Value = 114
$ws.Range("L34").Value = 7494
$ws.Range("M34").Value = -30
$ws.Range("N34").Value = -7662

$ws.Range("H39").Value = 1420
$ws.Range("I39").Value = 100
$ws.Range("J39").Value = 1750
$ws.Range("K39").Value = 300
$ws.Range("L39").Value = 5250
$ws.Range("M39").Value = -6
$ws.Range("N39").Value = -5838

$ws.Range("H55").Value = 1544.8334
$ws.Range("J55").Value = 1833
$ws.Range("L55").Value = 5499
$ws.Range("N55").Value = -5853

$ws.Range("H132").Value = 1176
$ws.Range("I132").Value = 981.125
$ws.Range("K132").Value = 8830.125
$ws.Range("M132").Value = -6300.125

$ws.Range("H140").Value = 3324.2354
$ws.Range("I140").Value = 3072.2856
$ws.Range("K140").Value = 9216.856800000001
$ws.Range("M140").Value = -4036.856800000001

$ws.Range("H141").Value = 89817.586
$ws.Range("I141").Value = 6777.8
$ws.Range("J141").Value = 505016.5
$ws.Range("K141").Value = 20333.4
$ws.Range("L141").Value = 1515049.5
$ws.Range("M141").Value = -15153.4
$ws.Range("N141").Value = -1525409.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10916.308
$ws.Range("I70").Value = 8997.25
$ws.Range("K70").Value = 8997.25
$ws.Range("M70").Value = -8727.25

$ws.Range("H73").Value = 10916.308
$ws.Range("I73").Value = 8997.25
$ws.Range("K73").Value = 8997.25
$ws.Range("M73").Value = -8061.25

$ws.Range("H101").Value = 49995
$ws.Range("J101").Value = 49995
$ws.Range("L101").Value = 49995
$ws.Range("N101").Value = -56485

$ws.Range("H104").Value = 30831.666
$ws.Range("J104").Value = 30831.666
$ws.Range("L104").Value = 30831.666
$ws.Range("N104").Value = -37819.666

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents() | Out-Null

$ws.Range("H123").Value = 35106.285
$ws.Range("J123").Value = 35106.285
$ws.Range("L123").Value = 35106.285
$ws.Range("N123").Value = -40006.285

$ws.Range("H126").Value = 29586
$ws.Range("I126").Value = 47933.145
$ws.Range("K126").Value = 143799.435
$ws.Range("M126").Value = -141329.435

$ws.Range("H132").Value = 838072.4399999999
$ws.Range("I132").Value = 913624.4399999999
$ws.Range("K132").Value = 2740873.32
$ws.Range("M132").Value = -2738343.32

$ws.Range("H134").Value = 50000.715
$ws.Range("J134").Value = 50000.715
$ws.Range("L134").Value = 150002.145
$ws.Range("N134").Value = -155072.145

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5066.769
$ws.Range("I40").Value = 4453.8
$ws.Range("K40").Value = 4453.8
$ws.Range("M40").Value = -4317.8

$ws.Range("H106").Value = 13118.182
$ws.Range("J106").Value = 13118.182
$ws.Range("L106").Value = 13118.182
$ws.Range("N106").Value = -15642.182

$ws.Range("H110").Value = 32500
$ws.Range("J110").Value = 32500
$ws.Range("L110").Value = 32500
$ws.Range("N110").Value = -40680

$ws.Range("H136").Value = 6760
$ws.Range("I136").Value = 6174.636
$ws.Range("K136").Value = 18523.908
$ws.Range("M136").Value = -15973.908

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 25500
$ws.Range("J40").Value = 25500
$ws.Range("L40").Value = 25500
$ws.Range("N40").Value = -25798

$ws.Range("H69").Value = 40853.668
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 40853.668
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 40853.668
$ws.Range("M69").ClearContents() | Out-Null
$ws.Range("N69").Value = -42351.668

$ws.Range("H72").Value = 40853.668
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 40853.668
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 122561.004
$ws.Range("M72").ClearContents() | Out-Null
$ws.Range("N72").Value = -130049.004

$ws.Range("H97").Value = 33060.25
$ws.Range("J97").Value = 33060.25
$ws.Range("L97").Value = 33060.25
$ws.Range("N97").Value = -35042.25

$ws.Range("H104").Value = 24183.5
$ws.Range("J104").Value = 24183.5
$ws.Range("L104").Value = 24183.5
$ws.Range("N104").Value = -31171.5

$ws.Range("H131").Value = 48425.43
$ws.Range("J131").Value = 48425.43
$ws.Range("L131").Value = 48425.43
$ws.Range("N131").Value = -58505.43

$ws.Range("H132").Value = 3177.8965
$ws.Range("I132").Value = 2948.261
$ws.Range("J132").Value = 4058.1667
$ws.Range("K132").Value = 8844.782999999999
$ws.Range("L132").Value = 12174.5001
$ws.Range("M132").Value = -6314.782999999999
$ws.Range("N132").Value = -17234.5001
